$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Surrogate Detection Statistics"

$win = $excel.ActiveWindow
$win.Width = 19218
$win.Height = 13344
$win.Left = 2952
$win.Top = 336

Write-Host "Width: $($win.Width)"
Write-Host "Height: $($win.Height)"
Write-Host "Left: $($win.Left)"
Write-Host "Top: $($win.Top)"
